# Updated cryptos list (price/volume refresh) on Tue Feb 28 18:38:15 UTC 2023
# with GitHub Actions.
#
# Note: Price ("D") and Volume(1h) ("E") cells hold plain text in this sheet
# (e.g. "1.000", "  +0.81%  "), not numbers. A leading "'" is used below for
# any new Price value that would otherwise be auto-parsed as a number, so it
# is stored/stays as text -- exactly like typing '1.000 directly into Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.488.62"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").Value = "1.641.86"
$ws.Range("E3").Value = "  +0.93%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.49%  "

$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").Value = "'303.46"
$ws.Range("E6").Value = "  +0.20%  "

$ws.Range("D7").Value = "'0.3805"
$ws.Range("E7").Value = "  +1.33%  "

$ws.Range("D8").Value = "'52.06"
$ws.Range("E8").Value = "  -0.76%  "

$ws.Range("D9").Value = "'0.3613"
$ws.Range("E9").Value = "  +0.30%  "

$ws.Range("D10").Value = "'0.08166"
$ws.Range("E10").Value = "  +1.43%  "

$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  -0.50%  "

$ws.Range("E13").Value = "  +0.02%  "

$ws.Range("D14").Value = "'6.452"
$ws.Range("E14").Value = "  -1.34%  "

$ws.Range("D15").Value = "'7.353"
$ws.Range("E15").Value = "  +2.22%  "

$ws.Range("D16").Value = "'0.00001237"
$ws.Range("E16").Value = "  -0.71%  "

$ws.Range("D17").Value = "1.635.94"
$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").Value = "'95.17"
$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("D19").Value = "'0.06956"
$ws.Range("E19").Value = "  +0.43%  "

$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'6.570"
$ws.Range("E20").Value = "  +1.94%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'17.53"
$ws.Range("E21").Value = "  -1.89%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("E23").Value = "  -1.21%  "

$ws.Range("D24").Value = "23.494.43"
$ws.Range("E24").Value = "  +0.77%  "

$ws.Range("D25").Value = "'2.524"
$ws.Range("E25").Value = "  +4.22%  "

$ws.Range("D26").Value = "'3.057"
$ws.Range("E26").Value = "  -4.40%  "

$ws.Range("E27").Value = "  +1.01%  "

$ws.Range("D28").Value = "'151.52"
$ws.Range("E28").Value = "  +2.39%  "

$ws.Range("D29").Value = "'5.272"
$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("D30").Value = "'133.26"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").Value = "1.817.65"
$ws.Range("E31").Value = "  +0.42%  "

$ws.Range("D32").Value = "'1.098"
$ws.Range("E32").Value = "  +16.09%  "

$ws.Range("D33").Value = "'2.153"
$ws.Range("E33").Value = "  -6.34%  "

$ws.Range("D34").Value = "'6.564"
$ws.Range("E34").Value = "  -2.65%  "

$ws.Range("D35").Value = "'11.49"
$ws.Range("E35").Value = "  +6.68%  "

$ws.Range("D36").Value = "'0.02768"
$ws.Range("E36").Value = "  -1.60%  "

$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("D38").Value = "'0.08754"
$ws.Range("E38").Value = "  -0.47%  "

$ws.Range("D39").Value = "'5.977"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("D40").Value = "'0.07029"
$ws.Range("E40").Value = "  -1.31%  "

$ws.Range("D41").Value = "'1.351"
$ws.Range("E41").Value = "  -0.78%  "

$ws.Range("D42").Value = "'0.7036"
$ws.Range("E42").Value = "  +0.41%  "

$ws.Range("D43").Value = "'12.26"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("D44").Value = "'15.64"
$ws.Range("E44").Value = "  -3.19%  "

$ws.Range("D45").Value = "'0.6531"
$ws.Range("E45").Value = "  +1.63%  "

$ws.Range("D46").Value = "'0.9998"
$ws.Range("E46").Value = "  -0.11%  "

$ws.Range("D47").Value = "'2.287"
$ws.Range("E47").Value = "  -0.80%  "

$ws.Range("D48").Value = "'3.963"
$ws.Range("E48").Value = "  -0.36%  "

$ws.Range("D49").Value = "'0.07975"
$ws.Range("E49").Value = "  -0.10%  "

$ws.Range("D50").Value = "'129.12"
$ws.Range("E50").Value = "  +2.30%  "

$ws.Range("D51").Value = "'1.191"
$ws.Range("E51").Value = "  -0.83%  "
